$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Fill in the "green cells" (algorithm results) for instance01 (row 4) and instance06 (row 9)
$ws.Range("I4").Value = 180.97507999999999
$ws.Range("I9").Value = 4.9803796

# Update the selected cell to reflect where the user ended up working
$ws.Range("K13").Select()
